$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "2023/24"
$ws.Range("B2").Value = "No"
$ws.Range("C2").Value = "Registered Trapline"
$ws.Range("D2").Value = "TR0615T010"
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = ""
$ws.Range("G2").Formula = "=""615"""
$ws.Range("G2").Copy()
$ws.Range("G2").PasteSpecial(-4163)
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = "No"
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = ""

$ws.Range("A2:M2").WrapText = $true
